$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all-classes")

# Row 3 - files/chessdoctor.pgn
$ws.Range("D3").Value = 33
$ws.Range("E3").Value = 45
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 26
$ws.Range("H3").Value = 12
$ws.Range("I3").Value = 32

# Row 5 - files/electronic_campfire.pgn
$ws.Range("D5").Value = 5
$ws.Range("F5").Value = 8
$ws.Range("H5").Value = 1

# Row 9 - files/GM_games.pgn
$ws.Range("D9").Value = 71
$ws.Range("E9").Value = 11
$ws.Range("H9").Value = 13
$ws.Range("I9").Value = 17

# Row 11 - files/hartwig.pgn
$ws.Range("D11").Value = 93
$ws.Range("E11").Value = 47
$ws.Range("F11").Value = 14
$ws.Range("G11").Value = 10
$ws.Range("H11").Value = 46
$ws.Range("I11").Value = 19

# Row 12 - files/hayes.pgn
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 1
$ws.Range("G12").Value = 1

# Row 13 - files/human_computer.pgn
$ws.Range("D13").Value = 55
$ws.Range("E13").Value = 19
$ws.Range("F13").Value = 1
$ws.Range("H13").Value = 11
$ws.Range("I13").Value = 15

# Row 14 - files/immortal_games.pgn
$ws.Range("E14").Value = 31
$ws.Range("G14").Value = 1

$excel.CalculateFullRebuild()
